$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 3277.9512
$ws.Range("J17").Value = 3277.9512
$ws.Range("L17").Value = 9833.8536
$ws.Range("N17").Value = -10169.8536
# Row 29
$ws.Range("H29").Value = 1124
$ws.Range("I29").Value = 655
$ws.Range("K29").Value = 1965
$ws.Range("M29").Value = -1684
# Row 38
$ws.Range("H38").Value = 561.6667
$ws.Range("I38").Value = 264.75
$ws.Range("J38").Value = 1155.5
$ws.Range("K38").Value = 794.25
$ws.Range("L38").Value = 3466.5
$ws.Range("M38").Value = -422.25
$ws.Range("N38").Value = -4210.5
# Row 58
$ws.Range("H58").Value = 37179.355
$ws.Range("I58").Value = 900
$ws.Range("J58").Value = 45066.176
$ws.Range("K58").Value = 2700
$ws.Range("L58").Value = 135198.528
$ws.Range("M58").Value = -2550
$ws.Range("N58").Value = -135498.528
# Row 87
$ws.Range("H87").Value = 39800
$ws.Range("J87").Value = 39800
$ws.Range("L87").Value = 39800
$ws.Range("N87").Value = -42296
# Row 90
$ws.Range("H90").Value = 39800
$ws.Range("J90").Value = 39800
$ws.Range("L90").Value = 119400
$ws.Range("N90").Value = -131880
# Row 132
$ws.Range("H132").Value = 2062.4443
$ws.Range("I132").Value = 1732.3667
$ws.Range("K132").Value = 5197.1001
$ws.Range("M132").Value = -2667.1001

$ws = $wb.Worksheets.Item("ARM")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
# Row 36
$ws.Range("H36").Value = 13378.25
$ws.Range("I36").Value = 8513
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 8513
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = -8167
$ws.Range("N36").Value = -15692
# Row 74
$ws.Range("H74").Value = 6160.1333
$ws.Range("I74").Value = 3505.8
$ws.Range("K74").Value = 3505.8
$ws.Range("M74").Value = -2631.8
# Row 77
$ws.Range("H77").Value = 6160.1333
$ws.Range("I77").Value = 3505.8
$ws.Range("K77").Value = 17529
$ws.Range("M77").Value = -13161
# Row 122
$ws.Range("H122").Value = 2485.5217
$ws.Range("I122").Value = 3047.7273
$ws.Range("J122").Value = 1970.1666
$ws.Range("K122").Value = 9143.1819
$ws.Range("L122").Value = 5910.4998
$ws.Range("M122").Value = -6693.1819
$ws.Range("N122").Value = -10810.4998

$ws = $wb.Worksheets.Item("BSM")
# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
# Row 119
$ws.Range("H119").Value = 67500
$ws.Range("J119").Value = 67500
$ws.Range("L119").Value = 67500
$ws.Range("N119").Value = -77176
# Row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
# Row 121
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 134
$ws.Range("H134").Value = 5657.037
$ws.Range("I134").Value = 5130.8335
$ws.Range("K134").Value = 15392.5005
$ws.Range("M134").Value = -12857.5005

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1282.4615
$ws.Range("I107").Value = 1282.4615
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1282.4615
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 637.5385000000001
$ws.Range("N107").ClearContents()
# Row 122
$ws.Range("H122").Value = 10856.154
$ws.Range("I122").Value = 3170.0588
$ws.Range("J122").Value = 25374.334
$ws.Range("K122").Value = 9510.1764
$ws.Range("L122").Value = 76123.00199999999
$ws.Range("M122").Value = -7060.1764
$ws.Range("N122").Value = -81023.00199999999
# Row 134
$ws.Range("H134").Value = 3065.4038
$ws.Range("I134").Value = 1769.8889
$ws.Range("J134").Value = 4464.56
$ws.Range("K134").Value = 5309.6667
$ws.Range("L134").Value = 13393.68
$ws.Range("M134").Value = -2774.6667
$ws.Range("N134").Value = -18463.68

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9331.25
$ws.Range("I80").Value = 18216.666
$ws.Range("K80").Value = 18216.666
$ws.Range("M80").Value = -17218.666
# Row 83
$ws.Range("H83").Value = 9331.25
$ws.Range("I83").Value = 18216.666
$ws.Range("K83").Value = 91083.33
$ws.Range("M83").Value = -86091.33
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 132
$ws.Range("H132").Value = 2959.25
$ws.Range("I132").Value = 2887.7144
$ws.Range("J132").Value = 3030.7856
$ws.Range("K132").Value = 8663.143199999999
$ws.Range("L132").Value = 9092.356800000001
$ws.Range("M132").Value = -6133.143199999999
$ws.Range("N132").Value = -14152.3568

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3959.1428
$ws.Range("I40").Value = 3584.3635
$ws.Range("K40").Value = 3584.3635
$ws.Range("M40").Value = -3448.3635
# Row 46
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -312
$ws.Range("N46").ClearContents()
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 103
$ws.Range("H103").Value = 60602
$ws.Range("J103").Value = 60602
$ws.Range("L103").Value = 60602
$ws.Range("N103").Value = -62946

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1903.174
$ws.Range("I132").Value = 1639.2
$ws.Range("J132").Value = 3663
$ws.Range("K132").Value = 4917.6
$ws.Range("L132").Value = 10989
$ws.Range("M132").Value = -2387.6
$ws.Range("N132").Value = -16049
# Row 136
$ws.Range("H136").Value = 3718.9155
$ws.Range("I136").Value = 1665.3182
$ws.Range("J136").Value = 7065.5186
$ws.Range("K136").Value = 4995.9546
$ws.Range("L136").Value = 21196.5558
$ws.Range("M136").Value = -2445.9546
$ws.Range("N136").Value = -26296.5558
